$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column values are plain text in the source data (e.g. "26.246.53" uses dots
# as thousands separators and would otherwise be auto-coerced to a number by Excel),
# so force the Text format before writing each one.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.246.53'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.660.01'
$ws.Range('E3').Value = '  -1.27%  '
$ws.Range('E4').Value = '  +0.54%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.48'
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5225'
$ws.Range('E6').Value = '  -1.65%  '
$ws.Range('E7').Value = '  +0.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2667'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06339'
$ws.Range('E9').Value = '  -1.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.27'
$ws.Range('E10').Value = '  -0.20%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07768'
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.445'
$ws.Range('E12').Value = '  -1.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.654.96'
$ws.Range('E13').Value = '  -1.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5498'
$ws.Range('E14').Value = '  -2.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0₅8303'
$ws.Range('E15').Value = '  -1.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.12'
$ws.Range('E16').Value = '  -1.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.258.77'
$ws.Range('E17').Value = '  -0.65%  '
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.700'
$ws.Range('E19').Value = '  -3.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '192.96'
$ws.Range('E20').Value = '  -1.26%  '
$ws.Range('E21').Value = '  -1.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.154'
$ws.Range('E22').Value = '  -3.89%  '
$ws.Range('E23').Value = '  +0.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '138.58'
$ws.Range('E24').Value = '  -3.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1242'
$ws.Range('E25').Value = '  -1.94%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.301'
$ws.Range('E26').Value = '  -2.80%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.17'
$ws.Range('E27').Value = '  -1.08%  '
$ws.Range('E28').Value = '  -1.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06056'
$ws.Range('E29').Value = '  -2.46%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.289'
$ws.Range('E30').Value = '  +0.94%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.567'
$ws.Range('E31').Value = '  +0.30%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.363'
$ws.Range('E32').Value = '  -2.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.655'
$ws.Range('E33').Value = '  -2.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9858'
$ws.Range('E34').Value = '  -3.44%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.783'
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5970'
$ws.Range('E37').Value = '  +3.72%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01599'
$ws.Range('E38').Value = '  -2.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.979'
$ws.Range('E39').Value = '  +0.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8667'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.048.20'
$ws.Range('E41').Value = '  -0.92%  '
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.11'
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.798.42'
$ws.Range('E44').Value = '  -1.68%  '
$ws.Range('E45').Value = '  +1.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.57'
$ws.Range('E46').Value = '  +0.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.009'
$ws.Range('E47').Value = '  +0.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.134'
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.486'
$ws.Range('E49').Value = '  +3.94%  '
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('E51').Value = '  +0.27%  '
